# The "reviews_count" column (column E) was removed from the sheet,
# shifting every later column (reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date) one slot to the
# left. Deleting the whole column reproduces that shift (and the updated
# dimension) in one step.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(5).Delete()
